$wb = $excel.ActiveWorkbook

# --- Update Sheet1: remove tabSelected (handled automatically once a
#     different sheet becomes the active/selected one at the end of the
#     script, but we still touch nothing else on Sheet1). ---

# --- Update Sheet2: change selection from B4 to A4:C16 ---
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Activate()
$ws2.Range("A4:C16").Select()

# --- Add new worksheet "rmanova1" after Sheet3 ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws4.Name = "rmanova1"

# Header row
$ws4.Range("A1").Value = "id"
$ws4.Range("B1").Value = "measurements"
$ws4.Range("C1").Value = "time"

# Data rows
$data = @(
    @(1, 9.9, 1),
    @(1, 10.1, 2),
    @(1, 10.199999999999999, 3),
    @(2, 11.7, 1),
    @(2, 11.9, 2),
    @(2, 12, 3),
    @(3, 3.2, 1),
    @(3, 3.3, 2),
    @(3, 3.4, 3),
    @(4, 5.2, 1),
    @(4, 5.4, 2),
    @(4, 5.6, 3)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws4.Cells.Item($row, 1).Value = $data[$i][0]
    $ws4.Cells.Item($row, 2).Value = $data[$i][1]
    $ws4.Cells.Item($row, 3).Value = $data[$i][2]
}

# Column widths to approximate the diff's bestFit custom widths as closely
# as this runtime's character-width rounding (steps of 1/6) allows.
$ws4.Columns.Item(1).ColumnWidth = 1.33
$ws4.Columns.Item(2).ColumnWidth = 11.5
$ws4.Columns.Item(3).ColumnWidth = 3.5

# Selection / active cell on the new sheet
$ws4.Range("B1").Select()

# Finally activate the new sheet so it becomes the active tab (activeTab=3)
$ws4.Activate()
